$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "même pas vrai!!!"
$ws.Range("D4").Value = ";)"

$ws.Range("A5").Select()
